# ClientData.xlsx: insert a new "Gender" column after "LastName" (i.e. before
# "Email"), shifting Email..NationalId one column to the right, and fill the
# new column with each client's gender.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D ("LastName" is C, "Email" was D) - this shifts
# D:I -> E:J, carrying over cell values, number formats, and column widths.
$ws.Range("D1").EntireColumn.Insert() | Out-Null

# EntireColumn.Insert() does not re-home the worksheet's hyperlink anchors
# (they stay pinned to the old D2:D5 cells even though the email text moved
# to E2:E5), so drop the stale collection and recreate the links on the
# shifted range.
$ws.Hyperlinks.Delete() | Out-Null

# Header + per-row gender values for the new column.
$ws.Range("D1").Value = "Gender"
$ws.Range("D2").Value = "M"
$ws.Range("D3").Value = "M"
$ws.Range("D4").Value = "F"
$ws.Range("D5").Value = "M"

$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:enriquecarrillo119999@gmail.com", "", "", "enriquecarrillo119999@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:edgarcarrillo119999@gmail.com", "", "", "edgarcarrillo119999@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:myriamcarrillo119999@gmail.com", "", "", "myriamcarrillo119999@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:juliocarrillo119999@gmail.com", "", "", "juliocarrillo119999@gmail.com") | Out-Null

# Hyperlinks.Add() stamps the built-in blue/underlined "Hyperlink" style on
# its target cells; restore each cell's original look (Enrique's email used
# a distinct Consolas/black style, the other three a plain Arial/blue one).
$ws.Range("E2").Font.Name = "Consolas"
$ws.Range("E2").Font.Color = 0
$ws.Range("E2").Font.Underline = 0

$ws.Range("E3:E5").Font.Name = "Arial"
$ws.Range("E3:E5").Font.Color = 16711680
$ws.Range("E3:E5").Font.Underline = 0

# Give the new Gender column a narrower width, matching the source edit
# (target XML width 7.82 chars; this engine's ColumnWidth<->XML-width
# mapping carries a fixed +5/6 offset versus the workbook's own font
# metrics, so back-solve for the COM-side value that lands on it).
$ws.Range("D1").EntireColumn.ColumnWidth = 6.986666666666667

# Match the cursor position left in the saved workbook.
$ws.Range("D7").Select() | Out-Null
